{"js": "// Update the date heading and every arithmetic-problem cell in the table,\n// in document order, matching the authored diff:\n//   \"2023-01-15 Sunday\" -> \"2023-01-16 Monday\"\n//   table cell 0..99 (row-major, 20 rows x 5 cols): old \"a\u00b1b=\" -> new \"a\u00b1b=\"\n\nconst OLD_TITLE = \"2023-01-15 Sunday\";\nconst NEW_TITLE = \"2023-01-16 Monday\";\n\n// Row-major (20 rows x 5 cols) replacement values for the table, in the\n// exact order the cells appear in the document.\nconst NEW_VALUES = [\n  [\"9+29=\", \"62+33=\", \"2+40=\", \"45+34=\", \"54+45=\"],\n  [\"21+10=\", \"70-8=\", \"82-64=\", \"67+21=\", \"64-30=\"],\n  [\"82-53=\", \"46-22=\", \"98-32=\", \"23-5=\", \"62-32=\"],\n  [\"71+25=\", \"67+29=\", \"85-49=\", \"12+60=\", \"12+51=\"],\n  [\"90-49=\", \"8+45=\", \"90-22=\", \"87-51=\", \"45-21=\"],\n  [\"86-59=\", \"77-37=\", \"89-79=\", \"77-23=\", \"58+20=\"],\n  [\"25+65=\", \"10+62=\", \"69-18=\", \"17+17=\", \"83-43=\"],\n  [\"24-21=\", \"55+21=\", \"88+1=\", \"77+20=\", \"85-78=\"],\n  [\"96-36=\", \"30+3=\", \"20+60=\", \"24+62=\", \"6+42=\"],\n  [\"38+43=\", \"75-27=\", \"85-48=\", \"57-57=\", \"54+41=\"],\n  [\"58-37=\", \"51-21=\", \"63+6=\", \"73-26=\", \"66-6=\"],\n  [\"58-4=\", \"52-13=\", \"63-43=\", \"69-2=\", \"98-33=\"],\n  [\"12+33=\", \"90-66=\", \"81-53=\", \"40+10=\", \"98-51=\"],\n  [\"19-11=\", \"50+45=\", \"81+17=\", \"55+7=\", \"14+85=\"],\n  [\"1+65=\", \"56+38=\", \"79+3=\", \"14+26=\", \"98-57=\"],\n  [\"42+4=\", \"77+11=\", \"54-35=\", \"59-5=\", \"98-26=\"],\n  [\"29+9=\", \"34+29=\", \"90-16=\", \"8+59=\", \"80-49=\"],\n  [\"40-18=\", \"6+26=\", \"33+57=\", \"23+31=\", \"84-41=\"],\n  [\"7-0=\", \"55+9=\", \"13+6=\", \"63-5=\", \"9+18=\"],\n  [\"9+70=\", \"18-2=\", \"79-51=\", \"16+68=\", \"85-12=\"],\n];\n\n// 1) Update the title paragraph (first paragraph of the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length > 0 && paragraphs.items[0].text === OLD_TITLE) {\n  paragraphs.items[0].getRange().insertText(NEW_TITLE, \"Replace\");\n} else {\n  // Fall back to a body-wide search/replace if the structure differs.\n  const hits = body.search(OLD_TITLE, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  for (const hit of hits.items) {\n    hit.insertText(NEW_TITLE, \"Replace\");\n  }\n}\nawait context.sync();\n\n// 2) Update every cell of the first table, in row-major document order.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length && r < NEW_VALUES.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length && r < NEW_VALUES.length; r++) {\n  const cells = rows.items[r].cells;\n  const rowVals = NEW_VALUES[r];\n  for (let c = 0; c < cells.items.length && c < rowVals.length; c++) {\n    cells.items[c].value = rowVals[c];\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date heading and every arithmetic-problem cell in the\n# table, in document order, matching the authored diff.\n\n$oldTitle = \"2023-01-15 Sunday\"\n$newTitle = \"2023-01-16 Monday\"\n\n$newValues = @(\n    @(\"9+29=\", \"62+33=\", \"2+40=\", \"45+34=\", \"54+45=\"),\n    @(\"21+10=\", \"70-8=\", \"82-64=\", \"67+21=\", \"64-30=\"),\n    @(\"82-53=\", \"46-22=\", \"98-32=\", \"23-5=\", \"62-32=\"),\n    @(\"71+25=\", \"67+29=\", \"85-49=\", \"12+60=\", \"12+51=\"),\n    @(\"90-49=\", \"8+45=\", \"90-22=\", \"87-51=\", \"45-21=\"),\n    @(\"86-59=\", \"77-37=\", \"89-79=\", \"77-23=\", \"58+20=\"),\n    @(\"25+65=\", \"10+62=\", \"69-18=\", \"17+17=\", \"83-43=\"),\n    @(\"24-21=\", \"55+21=\", \"88+1=\", \"77+20=\", \"85-78=\"),\n    @(\"96-36=\", \"30+3=\", \"20+60=\", \"24+62=\", \"6+42=\"),\n    @(\"38+43=\", \"75-27=\", \"85-48=\", \"57-57=\", \"54+41=\"),\n    @(\"58-37=\", \"51-21=\", \"63+6=\", \"73-26=\", \"66-6=\"),\n    @(\"58-4=\", \"52-13=\", \"63-43=\", \"69-2=\", \"98-33=\"),\n    @(\"12+33=\", \"90-66=\", \"81-53=\", \"40+10=\", \"98-51=\"),\n    @(\"19-11=\", \"50+45=\", \"81+17=\", \"55+7=\", \"14+85=\"),\n    @(\"1+65=\", \"56+38=\", \"79+3=\", \"14+26=\", \"98-57=\"),\n    @(\"42+4=\", \"77+11=\", \"54-35=\", \"59-5=\", \"98-26=\"),\n    @(\"29+9=\", \"34+29=\", \"90-16=\", \"8+59=\", \"80-49=\"),\n    @(\"40-18=\", \"6+26=\", \"33+57=\", \"23+31=\", \"84-41=\"),\n    @(\"7-0=\", \"55+9=\", \"13+6=\", \"63-5=\", \"9+18=\"),\n    @(\"9+70=\", \"18-2=\", \"79-51=\", \"16+68=\", \"85-12=\"),\n)\n\n$d = $word.ActiveDocument\n\n# 1) Update the title paragraph (first paragraph of the body).\n$titlePara = $d.Paragraphs(1).Range\nif ($titlePara.Text -eq $oldTitle) {\n    $titlePara.Text = $newTitle\n} else {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldTitle\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newTitle\n    $find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)\n}\n\n# 2) Update every cell of the first table, in row-major document order.\n$t = $d.Tables(1)\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n    }\n}\n"}
